$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete rows that were removed entirely (original row numbers) ---
# Row 26 "RM 232" and row 28 "SC 92" (original numbering) were removed,
# shifting everything below up.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# --- Update individual cell values (using final row numbers after deletion) ---
$ws.Range("E5").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("E11").Value = -7.9
$ws.Range("D19").Value = -15.5
$ws.Range("E19").Value = ""
$ws.Range("D21").Value = ""
$ws.Range("D23").Value = -13.9
$ws.Range("E23").Value = -7
$ws.Range("F24").Value = 16.78
$ws.Range("E25").Value = -7.1

$ws.Range("B26").Value = ""
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = ""
$ws.Range("F28").Value = 17.44
$ws.Range("B29").Value = ""
$ws.Range("E29").Value = ""
$ws.Range("F30").Value = ""
$ws.Range("F32").Value = ""
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
